$wb = $excel.ActiveWorkbook

# --- donationShortDetail sheet: update donation code + requestId, move selection ---
$ws = $wb.Worksheets.Item("donationShortDetail")
$ws.Activate()
$ws.Range("A3").Value = "CDF33"
$ws.Range("B3").Value = "CDF33-20200521"
$ws.Range("E28").Select()

# --- orderTests sheet: update donation code + requestId, move selection ---
$ws = $wb.Worksheets.Item("orderTests")
$ws.Activate()
$ws.Range("A3").Value = "CDF33"
$ws.Range("B3").Value = "CDF33-20200521"
$ws.Range("B3").Select()

# --- orderTestsStatusHistory sheet: update donation code + requestId, move selection ---
$ws = $wb.Worksheets.Item("orderTestsStatusHistory")
$ws.Activate()
$ws.Range("A3").Value = "CDF33"
$ws.Range("B3").Value = "CDF33-20200521"
$ws.Range("E19").Select()
